# Atualiza instrução de trabalho
# - Marca como "Resolvido" um conjunto de incidentes que estavam "Pendente"
#   (aba SPN e aba ITI, coluna Status).
# - Corrige o nome do responsável "Erick Silva" para "Erick da Silva" em
#   várias linhas da aba ITI (coluna Responsavel).

$wb = $excel.ActiveWorkbook

# --- Aba SPN: linhas cuja Status passa de "Pendente" para "Resolvido" ---
$wsSPN = $wb.Worksheets.Item("SPN")
$spnResolvidoRows = @(104, 106, 107, 109)
foreach ($r in $spnResolvidoRows) {
    $wsSPN.Range("J$r").Value = "Resolvido"
}

# --- Aba ITI: linhas cuja Status passa de "Pendente" para "Resolvido" ---
$wsITI = $wb.Worksheets.Item("ITI")
$itiResolvidoRows = @(312, 317, 325, 327, 331, 332, 335, 336, 337, 338, 339, 340, 344, 345, 346, 347)
foreach ($r in $itiResolvidoRows) {
    $wsITI.Range("J$r").Value = "Resolvido"
}

# --- Aba ITI: linhas cujo Responsavel "Erick Silva" passa a "Erick da Silva" ---
$itiErickRows = @(
    116, 117, 118,
    125, 126, 127, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141, 142, 143, 144, 145,
    153,
    161, 162, 163, 164, 165, 166, 167, 168,
    180, 181, 182, 183, 184, 185, 186, 187, 188, 189,
    234, 235,
    245
)
foreach ($r in $itiErickRows) {
    $wsITI.Range("B$r").Value = "Erick da Silva"
}
